$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022317612170957
$ws.Range("D2").Value = 1.027891340981427
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.033169181054141
$ws.Range("I2").Value = 1.030868811031093
$ws.Range("J2").Value = 1.027504127111424
$ws.Range("K2").Value = 1.030709749231591
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.035972304702392
$ws.Range("N2").Value = 1.013166432372473
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023580389629369
$ws.Range("D3").Value = 1.028854726943839
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.034606232269697
$ws.Range("I3").Value = 1.031192233512502
$ws.Range("J3").Value = 1.028403274761367
$ws.Range("K3").Value = 1.03148053557865
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.037216603349391
$ws.Range("N3").Value = 1.013468901000798
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024396200586253
$ws.Range("D4").Value = 1.029476698461201
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.035534940328097
$ws.Range("I4").Value = 1.031399180480371
$ws.Range("J4").Value = 1.028983356342545
$ws.Range("K4").Value = 1.031977272631944
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.038020081683448
$ws.Range("N4").Value = 1.013663897200889
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0247388629429
$ws.Range("D5").Value = 1.029737841705261
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.035925096817883
$ws.Range("I5").Value = 1.031485624602773
$ws.Range("J5").Value = 1.029226812132671
$ws.Range("K5").Value = 1.032185621261896
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.03835746980496
$ws.Range("N5").Value = 1.013745701923837
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02479637972696
$ws.Range("D6").Value = 1.029781669314242
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.035990589992313
$ws.Range("I6").Value = 1.031500106358869
$ws.Range("J6").Value = 1.029267665444275
$ws.Range("K6").Value = 1.032220575855875
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.038414095776678
$ws.Range("N6").Value = 1.013759427257862
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024400780448774
$ws.Range("D7").Value = 1.029480189177762
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.035540154681837
$ws.Range("I7").Value = 1.031400337735732
$ws.Range("J7").Value = 1.028986611020103
$ws.Range("K7").Value = 1.031980058478269
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.038024591423072
$ws.Range("N7").Value = 1.01366499095297
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022744643575323
$ws.Range("D8").Value = 1.028217213888345
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.033655082703521
$ws.Range("I8").Value = 1.030978596127215
$ws.Range("J8").Value = 1.027808357659325
$ws.Range("K8").Value = 1.030970658366101
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.036393169068021
$ws.Range("N8").Value = 1.013268802990971
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019816237852612
$ws.Range("D9").Value = 1.025980828910236
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.030324217739091
$ws.Range("I9").Value = 1.030217542103066
$ws.Range("J9").Value = 1.025718763693928
$ws.Range("K9").Value = 1.029176440694046
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.033505422497093
$ws.Range("N9").Value = 1.012565101012125
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01785691072253
$ws.Range("D10").Value = 1.024482444196289
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.028097163381579
$ws.Range("I10").Value = 1.029698065791922
$ws.Range("J10").Value = 1.02431652801195
$ws.Range("K10").Value = 1.027969704745342
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.031571237540847
$ws.Range("N10").Value = 1.012092159550601
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017006763784151
$ws.Range("D11").Value = 1.023831820860453
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.027131211053105
$ws.Range("I11").Value = 1.029470237211614
$ws.Range("J11").Value = 1.023707125627495
$ws.Range("K11").Value = 1.027444627910202
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.030731507103183
$ws.Range("N11").Value = 1.011886453417033
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01669071360931
$ws.Range("D12").Value = 1.023589874755975
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.02677216277606
$ws.Range("I12").Value = 1.029385175414259
$ws.Range("J12").Value = 1.023480428330496
$ws.Range("K12").Value = 1.0272492044352
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.030419255039934
$ws.Range("N12").Value = 1.011809905602325
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016758519650113
$ws.Range("D13").Value = 1.023641785542289
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.026849191297754
$ws.Range("I13").Value = 1.029403441220402
$ws.Range("J13").Value = 1.023529071032568
$ws.Range("K13").Value = 1.027291140977036
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.030486249594023
$ws.Range("N13").Value = 1.011826331697038
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01698064447224
$ws.Range("D14").Value = 1.023811827161742
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.027101537161205
$ws.Range("I14").Value = 1.029463214883167
$ws.Range("J14").Value = 1.023688393677955
$ws.Range("K14").Value = 1.027428482048655
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.030705703195868
$ws.Range("N14").Value = 1.011880128801868
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017117467260714
$ws.Range("D15").Value = 1.023916558793785
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.02725698237078
$ws.Range("I15").Value = 1.029499985570167
$ws.Range("J15").Value = 1.023786512699187
$ws.Range("K15").Value = 1.027513051098943
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.030840870731915
$ws.Range("N15").Value = 1.011913256463842
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017913294931304
$ws.Range("D16").Value = 1.024525585451378
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.028161235687069
$ws.Range("I16").Value = 1.029713124955578
$ws.Range("J16").Value = 1.024356924835655
$ws.Range("K16").Value = 1.028004498328068
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.031626920500094
$ws.Range("N16").Value = 1.012105792110838
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.01841202543112
$ws.Range("D17").Value = 1.024907124133014
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.028728010624464
$ws.Range("I17").Value = 1.0298460462618
$ws.Range("J17").Value = 1.024714130892126
$ws.Range("K17").Value = 1.028312084701121
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.032119391410217
$ws.Range("N17").Value = 1.012226317600404
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018702758468349
$ws.Range("D18").Value = 1.025129494499088
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.02905844470888
$ws.Range("I18").Value = 1.029923298096074
$ws.Range("J18").Value = 1.024922268714411
$ws.Range("K18").Value = 1.028491248591042
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.032406427991508
$ws.Range("N18").Value = 1.012296529461377
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018801862587002
$ws.Range("D19").Value = 1.025205287514431
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.029171087909319
$ws.Range("I19").Value = 1.029949591715549
$ws.Range("J19").Value = 1.024993202077335
$ws.Range("K19").Value = 1.028552297208772
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.032504264035086
$ws.Range("N19").Value = 1.012320454898654
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018358533789585
$ws.Range("D20").Value = 1.02486620673725
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.028667217238998
$ws.Range("I20").Value = 1.029831813925514
$ws.Range("J20").Value = 1.024675828270825
$ws.Range("K20").Value = 1.02827910906665
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.032066576056383
$ws.Range("N20").Value = 1.012213395531796
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016915241684241
$ws.Range("D21").Value = 1.023761761778095
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.027027234565264
$ws.Range("I21").Value = 1.029445625081596
$ws.Range("J21").Value = 1.023641486494496
$ws.Range("K21").Value = 1.027388049222348
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.030641088979921
$ws.Range("N21").Value = 1.011864290753697
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016006234770549
$ws.Range("D22").Value = 1.023065757051517
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.025994659261596
$ws.Range("I22").Value = 1.029200288752341
$ws.Range("J22").Value = 1.022989196602336
$ws.Range("K22").Value = 1.02682556644975
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.029742865576927
$ws.Range("N22").Value = 1.011643987573349
$ws.Range("B23").Value = 1.019999999999999
$ws.Range("C23").Value = 1.016488265491374
$ws.Range("D23").Value = 1.023434874740305
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.02654218696172
$ws.Range("I23").Value = 1.029330585988735
$ws.Range("J23").Value = 1.023335174712747
$ws.Range("K23").Value = 1.027123962404061
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.030219218866317
$ws.Range("N23").Value = 1.011760851394056
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018382704878692
$ws.Range("D24").Value = 1.024884696086747
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.028694687635071
$ws.Range("I24").Value = 1.02983824576737
$ws.Range("J24").Value = 1.02469313623899
$ws.Range("K24").Value = 1.028294010096493
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.032090441700233
$ws.Range("N24").Value = 1.012219234732121
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020574524134733
$ws.Range("D25").Value = 1.026560291378501
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.031186442751063
$ws.Range("I25").Value = 1.030416420972698
$ws.Range("J25").Value = 1.02626057659954
$ws.Range("K25").Value = 1.029642144403767
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.034253542073681
$ws.Range("N25").Value = 1.012747691158529
